# v1.1 update Owner Status
# Change the "Owner Status" column (I) values from "Open" to "Closed"
# for the three existing review rows on the LH_TC_PUBLISH&UPLOAD_REVIEWS sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LH_TC_PUBLISH&UPLOAD_REVIEWS")

$ws.Range("I2").Value = "Closed"
$ws.Range("I3").Value = "Closed"
$ws.Range("I4").Value = "Closed"
